# "Generate Report for handoff"
#
# Source file 43073d14-65d6-4022-a602-b675bf8408ea.md finished its
# handoff/handback cycle and is replaced in the report by a new source
# file dff677dc-47bd-4419-a55b-4d35806f3fa6.md (whose handoff artifacts
# use transform hash 990e10d5bb313e844bffddbc64e112d9f5777277 and new
# handoff timestamps). In addition a second source file,
# f5864e5f-6d52-4c13-892b-d6126c317490.md, shows up whose handoff
# transform failed, so it is listed with an "Ignored" row (no handoff
# file/datetime yet) ahead of the ".localization-config" housekeeping row.

$wb = $excel.ActiveWorkbook

$oldUuid  = "43073d14-65d6-4022-a602-b675bf8408ea"
$newUuid  = "dff677dc-47bd-4419-a55b-4d35806f3fa6"
$failUuid = "f5864e5f-6d52-4c13-892b-d6126c317490"
$oldHash  = "4b286addc56cd24433acde6cb8c8c7ff7bf179f7"
$newHash  = "990e10d5bb313e844bffddbc64e112d9f5777277"

$newMdName   = "$newUuid.md"
$failMdName  = "$failUuid.md"
$cfgName     = ".localization-config"

$zhXlfName = "$newUuid.$newHash.zh-cn.xlf"
$deXlfName = "$newUuid.$newHash.de-de.xlf"

$zhHandoffDt = "2016-02-15 08:38:16"
$deHandoffDt = "2016-02-15 08:38:30"
$epochDt     = "0001-01-01 00:00:00"

$repoBase    = "https://github.com/OpenLocalizationTest/oltest/blob/72e352ba8404bfc5802eb5384cf496161d0b61a9"
$zhHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f0459038964029d15126abca61b32ea4a404f163/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht"
$deHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9c86082c17355e5983265ec430b1287293d904df/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht"

function Set-FileNameLink($ws, [string]$cellRef, [string]$display, [string]$url) {
    $ws.Range($cellRef).Value = $display
    $ws.Range($cellRef).Style = "HyperLink"
    [void]$ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $display)
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()

$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

$ws1.Range("B3").Value = "Handoff transform failed"
$ws1.Range("C3").Value = "Handoff transform failed"

$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

Set-FileNameLink $ws1 "A2" $newMdName  "$repoBase/e2e/$newMdName"
Set-FileNameLink $ws1 "A3" $failMdName "$repoBase/e2e/$failMdName"
Set-FileNameLink $ws1 "A4" $cfgName    "$repoBase/$cfgName"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("D2").Value = $zhHandoffDt
$ws2.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G2").Value = $epochDt
$ws2.Range("H2").Value = "Include"

$ws2.Range("B3").Value = "Handoff transform failed"
$ws2.Range("D3").Value = $epochDt
$ws2.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G3").Value = $epochDt
$ws2.Range("H3").Value = "Ignored"

$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = $epochDt
$ws2.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G4").Value = $epochDt
$ws2.Range("H4").Value = "Ignored"

Set-FileNameLink $ws2 "A2" $newMdName  "$repoBase/e2e/$newMdName"
Set-FileNameLink $ws2 "C2" $zhXlfName  "$zhHandoffBase/$zhXlfName"
Set-FileNameLink $ws2 "A3" $failMdName "$repoBase/e2e/$failMdName"
Set-FileNameLink $ws2 "A4" $cfgName    "$repoBase/$cfgName"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("D2").Value = $deHandoffDt
$ws3.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G2").Value = $epochDt
$ws3.Range("H2").Value = "Include"

$ws3.Range("B3").Value = "Handoff transform failed"
$ws3.Range("D3").Value = $epochDt
$ws3.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G3").Value = $epochDt
$ws3.Range("H3").Value = "Ignored"

$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = $epochDt
$ws3.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G4").Value = $epochDt
$ws3.Range("H4").Value = "Ignored"

Set-FileNameLink $ws3 "A2" $newMdName  "$repoBase/e2e/$newMdName"
Set-FileNameLink $ws3 "C2" $deXlfName  "$deHandoffBase/$deXlfName"
Set-FileNameLink $ws3 "A3" $failMdName "$repoBase/e2e/$failMdName"
Set-FileNameLink $ws3 "A4" $cfgName    "$repoBase/$cfgName"
